# Drop in RMI script files
# - Remove the "Texas Notes" sheet (it was a scratch/notes tab, no longer needed).
# - Fill in the real sales-volume numbers for Gas Boilers / Oil Boilers / Oil
#   Furnaces on the Data tab (previously placeholder 0s highlighted yellow),
#   which also ripples into the ENERGY STAR fraction calcs on the three
#   BFoCSbQL-* tabs via their existing formulas.

$wb = $excel.ActiveWorkbook

# --- Delete the "Texas Notes" worksheet -------------------------------------
$notes = $wb.Worksheets.Item("Texas Notes")
$notes.Delete()

# --- Fill in the real Data values (was 0 placeholders w/ yellow highlight) --
$data = $wb.Worksheets.Item("Data")

$data.Range("C9").Value = 192000
$data.Range("C10").Value = 123000
$data.Range("C11").Value = 56000

# Clear the placeholder yellow highlight now that real data is in place.
$data.Range("C9:C10").Interior.ColorIndex = -4142
$data.Range("C9:C10").HorizontalAlignment = -4131

$data.Range("C11").Interior.ColorIndex = -4142
$data.Range("C11").HorizontalAlignment = -4131
$data.Range("C11").Borders.Item(9).LineStyle = 1

$wb.Application.Calculate()
